$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Fixed Excel outputs": corrected computed result values ---
# Row 10 (Conductivity / parameter row) — corrected from 0.045 to 0.001 for all materials
$ws.Range("B10").Value2 = 0.001
$ws.Range("C10").Value2 = 0.001
$ws.Range("D10").Value2 = 0.001

# Row 11 — corrected Bi for lumped capacitance results
$ws.Range("B11").Value2 = 0.24998417821656865
$ws.Range("C11").Value2 = 0.24993803189291913
$ws.Range("D11").Value2 = 0.24630541871921183

# Row 12 — corrected results
$ws.Range("B12").Value2 = 0.002573577621669515
$ws.Range("C12").Value2 = 0.014081781376518221
$ws.Range("D12").Value2 = 0.62541871921182257

# Row 18 — corrected steady-state flux results
$ws.Range("B18").Value2 = 2553.6666666666665
$ws.Range("C18").Value2 = 2553.6666666666665
$ws.Range("D18").Value2 = 2553.6666666666665

# Row 19 — corrected temperature results
$ws.Range("B19").Value2 = 126.98267125644384
$ws.Range("C19").Value2 = 126.98267125644384
$ws.Range("D19").Value2 = 126.98267125644384

# --- Formatting fix: remove the red highlight fill on B18:D19 (no longer flagged) ---
$ws.Range("B18:D19").Interior.Pattern = -4142

# --- Column width touch-ups (columns already had customWidth set; nudging closer
#     to the re-saved layout's measured widths) ---
$ws.Columns.Item(2).ColumnWidth = 15.0
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(4).ColumnWidth = 12.833333333333334
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.666666666666666

# --- Update the saved cursor/selection position ---
$ws.Range("E22").Select()
